$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match column width used across B:J for the newly-touched columns K:N
$ws.Range("K1:N1").ColumnWidth = 8.7109375

# Extend the table one more year (2023) into column K, mirroring column J's layout/format
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 376
$ws.Range("K5").Value = 22
$ws.Range("K6").Value = 354

# Close off the table with a right border along the new last column
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2
$ws.Range("K3:K6").Borders.Item(10).Color = 0
